$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$new = $wb.Worksheets.Add($null, $ws1)
$new.Name = "Sheet1"

$src = $ws1.Range("A9:K40")
$src.Copy()
$new.Range("A1").PasteSpecial(-4104)  # xlPasteAll

$header = $new.Range("A1:K1")
$header.HorizontalAlignment = -4108  # xlCenter
$header.VerticalAlignment = -4108    # xlCenter
$header.Borders.LineStyle = 1
$header.Borders.Weight = 2
$header.Borders.Color = 0

$body = $new.Range("A2:K32")
$body.HorizontalAlignment = -4131  # xlLeft
$body.VerticalAlignment = -4160    # xlTop
$body.Borders.LineStyle = 1
$body.Borders.Weight = 2
$body.Borders.Color = 0

Write-Output "ok"
